# Applies the APP-000029 assets_liabilities.xlsx edit:
#   - Summary sheet: new client name, updated income/assets/liabilities/net-worth/ratio
#   - Assets sheet: vehicle descriptions + values updated
#   - Liabilities sheet: auto loan figures updated and a new "Personal Loans" row inserted

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")

$summary.Range("B3").Value = "Mona Al Kalbani"
$summary.Range("B4").Value = 1259.85
$summary.Range("B6").Value = 316564
$summary.Range("B7").Value = 120533
$summary.Range("B8").Value = 196031
$summary.Range("B9").Value = 2.63

# ---------------------------------------------------------------------------
# Assets sheet
# ---------------------------------------------------------------------------
$assets = $wb.Worksheets.Item("Assets")

$assets.Range("B2").Value = "Premium Car"
$assets.Range("C2").Value = 204534
$assets.Range("B3").Value = "Mid-range Car"
$assets.Range("C3").Value = 109614
$assets.Range("C4").Value = 2416
$assets.Range("C5").Value = 316564

# ---------------------------------------------------------------------------
# Liabilities sheet
# ---------------------------------------------------------------------------
$liab = $wb.Worksheets.Item("Liabilities")

# Update the existing Auto Loans row (row 2)
$liab.Range("C2").Value = 65768
$liab.Range("D2").Value = 783

# Insert a new row for "Personal Loans" above the current Credit Cards row (row 3),
# pushing Credit Cards to row 4 and TOTAL LIABILITIES to row 5.
$liab.Rows.Item(3).Insert()

# Copy the formatting from row 2 (Auto Loans) onto the new blank row so the new
# row matches the sheet's existing look (border + number format), same as
# every other data row.
$liab.Range("A2:E2").Copy()
$liab.Range("A3:E3").PasteSpecial(-4122)

# Populate the new Personal Loans row
$liab.Range("A3").Value = "Personal Loans"
$liab.Range("B3").Value = "Personal Loan"
$liab.Range("C3").Value = 43472
$liab.Range("D3").Value = 725
$liab.Range("E3").Value = 5

# Update the (now shifted) Credit Cards row (row 4)
$liab.Range("C4").Value = 11293
$liab.Range("D4").Value = 565

# Update the (now shifted) TOTAL LIABILITIES row (row 5)
$liab.Range("C5").Value = 120533
